# Remove the obsolete "Annick" reservation row (row 37). Excel shifts all
# subsequent rows up by one, which reproduces the rest of the diff
# (every following row's content moving up one position) without having
# to touch each cell individually.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(37).Delete()
